$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

$ws.Range("E10").Value = 506
$ws.Range("F10").Value = 255
$ws.Range("H10").Value = 255
$ws.Range("E11").Value = 333
$ws.Range("F11").Value = 186
$ws.Range("H11").Value = 186
$ws.Range("E12").Value = 494
$ws.Range("F12").Value = 269
$ws.Range("H12").Value = 269
$ws.Range("F13").Value = 68
$ws.Range("H13").Value = 68
$ws.Range("E14").Value = 124
$ws.Range("F14").Value = 66
$ws.Range("H14").Value = 66
$ws.Range("F15").Value = 70
$ws.Range("H15").Value = 70
$ws.Range("F16").Value = 100
$ws.Range("H16").Value = 100
$ws.Range("F17").Value = 53
$ws.Range("H17").Value = 53
$ws.Range("F20").Value = 30
$ws.Range("H20").Value = 30
$ws.Range("E21").Value = 134
$ws.Range("F21").Value = 74
$ws.Range("H21").Value = 74
$ws.Range("F22").Value = 87
$ws.Range("H22").Value = 87
$ws.Range("F24").Value = 111
$ws.Range("H24").Value = 111
$ws.Range("F26").Value = 94
$ws.Range("H26").Value = 94
$ws.Range("E27").Value = 315
$ws.Range("F27").Value = 157
$ws.Range("H27").Value = 157
$ws.Range("F28").Value = 74
$ws.Range("H28").Value = 74
$ws.Range("F29").Value = 93
$ws.Range("H29").Value = 93
$ws.Range("E30").Value = 198
$ws.Range("F30").Value = 119
$ws.Range("H30").Value = 119
$ws.Range("E33").Value = 282
$ws.Range("F33").Value = 143
$ws.Range("H33").Value = 143
$ws.Range("F37").Value = 74
$ws.Range("H37").Value = 74
$ws.Range("F38").Value = 55
$ws.Range("H38").Value = 55
$ws.Range("E39").Value = 178
$ws.Range("F39").Value = 88
$ws.Range("H39").Value = 88
$ws.Range("E40").Value = 254
$ws.Range("F41").Value = 177
$ws.Range("H41").Value = 177
$ws.Range("E42").Value = 367
$ws.Range("F42").Value = 203
$ws.Range("H42").Value = 203
$ws.Range("E45").Value = 138
$ws.Range("F45").Value = 70
$ws.Range("H45").Value = 70
$ws.Range("E46").Value = 303
$ws.Range("E47").Value = 438
$ws.Range("F47").Value = 221
$ws.Range("H47").Value = 221
$ws.Range("E48").Value = 199
$ws.Range("F48").Value = 87
$ws.Range("H48").Value = 87
$ws.Range("E49").Value = 279
$ws.Range("F49").Value = 120
$ws.Range("H49").Value = 120
$ws.Range("E50").Value = 236
$ws.Range("F50").Value = 111
$ws.Range("H50").Value = 111
$ws.Range("F51").Value = 98
$ws.Range("H51").Value = 98
